$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.895.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.712.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.51%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +20.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "675.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.429"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.711.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.46%  "
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.402.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.669.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  +14.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.712.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("E20").Value = "  +5.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.553"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "518.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +7.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("E29").Value = "  +11.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("E31").Value = "  +7.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.187"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.598"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.89%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "616.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.164"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.971"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.32%  "
$ws.Range("E43").Value = "  +9.36%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0448"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.426"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +24.68%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.72%  "
